$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new column before column D. This shifts the existing
#    reporting-period columns right by one (old D->E, E->F, ... K->L)
#    and leaves a blank column D ready for the newest period.
# ------------------------------------------------------------------
$ws.Columns("D").Insert()

# ------------------------------------------------------------------
# 2) The newly inserted column does not inherit the row-specific
#    number formatting (date format row 7/38/80, integer format for
#    the rest). Copy the formatting from the neighboring column E
#    (which used to be D) across to the new column D.
# ------------------------------------------------------------------
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Match the column width of the new column D to the rest of the
#    data columns.
# ------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# ------------------------------------------------------------------
# 4) Write the refreshed figures: one new reporting period in column
#    D plus the restated historical figures for columns E:K. The
#    oldest period (now column K) keeps its original value.
# ------------------------------------------------------------------
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43100
$ws.Cells.Item(7, 6).Value = 42735
$ws.Cells.Item(7, 7).Value = 42369
$ws.Cells.Item(7, 8).Value = 42004
$ws.Cells.Item(7, 9).Value = 41639
$ws.Cells.Item(7, 10).Value = 41274
$ws.Cells.Item(7, 11).Value = 40908

$ws.Cells.Item(8, 4).Value = 1017500
$ws.Cells.Item(8, 5).Value = 1068000
$ws.Cells.Item(8, 6).Value = 1111000
$ws.Cells.Item(8, 7).Value = 1153300
$ws.Cells.Item(8, 8).Value = 1260900
$ws.Cells.Item(8, 9).Value = 1359400
$ws.Cells.Item(8, 10).Value = 1638300
$ws.Cells.Item(8, 11).Value = 1817300

$ws.Cells.Item(9, 4).Value = 734200
$ws.Cells.Item(9, 5).Value = 739400
$ws.Cells.Item(9, 6).Value = 745500
$ws.Cells.Item(9, 7).Value = 762300
$ws.Cells.Item(9, 8).Value = 752400
$ws.Cells.Item(9, 9).Value = 824900
$ws.Cells.Item(9, 10).Value = 955400
$ws.Cells.Item(9, 11).Value = 952000

$ws.Cells.Item(10, 4).Value = 283300
$ws.Cells.Item(10, 5).Value = 328600
$ws.Cells.Item(10, 6).Value = 365600
$ws.Cells.Item(10, 7).Value = 391000
$ws.Cells.Item(10, 8).Value = 508500
$ws.Cells.Item(10, 9).Value = 534400
$ws.Cells.Item(10, 10).Value = 682900
$ws.Cells.Item(10, 11).Value = 865400

$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(12, 6).Value = "NA"
$ws.Cells.Item(12, 7).Value = "NA"
$ws.Cells.Item(12, 8).Value = "NA"
$ws.Cells.Item(12, 9).Value = "NA"
$ws.Cells.Item(12, 10).Value = "NA"
$ws.Cells.Item(12, 11).Value = "NA"

$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0

$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0

$ws.Cells.Item(15, 4).Value = 53500
$ws.Cells.Item(15, 5).Value = 39500
$ws.Cells.Item(15, 6).Value = 38100
$ws.Cells.Item(15, 7).Value = 49900
$ws.Cells.Item(15, 8).Value = 55200
$ws.Cells.Item(15, 9).Value = 71200
$ws.Cells.Item(15, 10).Value = 91000
$ws.Cells.Item(15, 11).Value = 82700

$ws.Cells.Item(17, 4).Value = 997100
$ws.Cells.Item(17, 5).Value = 986100
$ws.Cells.Item(17, 6).Value = 1025500
$ws.Cells.Item(17, 7).Value = 1067700
$ws.Cells.Item(17, 8).Value = 1078200
$ws.Cells.Item(17, 9).Value = 1179700
$ws.Cells.Item(17, 10).Value = 1366500
$ws.Cells.Item(17, 11).Value = 1420100

$ws.Cells.Item(18, 4).Value = 20400
$ws.Cells.Item(18, 5).Value = 81900
$ws.Cells.Item(18, 6).Value = 85500
$ws.Cells.Item(18, 7).Value = 85500
$ws.Cells.Item(18, 8).Value = 182600
$ws.Cells.Item(18, 9).Value = 179600
$ws.Cells.Item(18, 10).Value = 271800
$ws.Cells.Item(18, 11).Value = 397200

$ws.Cells.Item(20, 4).Value = -1700
$ws.Cells.Item(20, 5).Value = 800
$ws.Cells.Item(20, 6).Value = 1900
$ws.Cells.Item(20, 7).Value = -2200
$ws.Cells.Item(20, 8).Value = 14600
$ws.Cells.Item(20, 9).Value = 17100
$ws.Cells.Item(20, 10).Value = 24000
$ws.Cells.Item(20, 11).Value = 32400

$ws.Cells.Item(21, 4).Value = 179500
$ws.Cells.Item(21, 5).Value = 235600
$ws.Cells.Item(21, 6).Value = 234500
$ws.Cells.Item(21, 7).Value = 238100
$ws.Cells.Item(21, 8).Value = 365200
$ws.Cells.Item(21, 9).Value = 382800
$ws.Cells.Item(21, 10).Value = 506400
$ws.Cells.Item(21, 11).Value = 636100

$ws.Cells.Item(22, 4).Value = 38100
$ws.Cells.Item(22, 5).Value = 40600
$ws.Cells.Item(22, 6).Value = 43300
$ws.Cells.Item(22, 7).Value = 46600
$ws.Cells.Item(22, 8).Value = 69300
$ws.Cells.Item(22, 9).Value = 85000
$ws.Cells.Item(22, 10).Value = 95500
$ws.Cells.Item(22, 11).Value = 114200

$ws.Cells.Item(23, 4).Value = -19300
$ws.Cells.Item(23, 5).Value = 42200
$ws.Cells.Item(23, 6).Value = 44100
$ws.Cells.Item(23, 7).Value = 36700
$ws.Cells.Item(23, 8).Value = 128000
$ws.Cells.Item(23, 9).Value = 111700
$ws.Cells.Item(23, 10).Value = 200300
$ws.Cells.Item(23, 11).Value = 315400

$ws.Cells.Item(24, 4).Value = -1700
$ws.Cells.Item(24, 5).Value = 11000
$ws.Cells.Item(24, 6).Value = 2800
$ws.Cells.Item(24, 7).Value = 9900
$ws.Cells.Item(24, 8).Value = 30300
$ws.Cells.Item(24, 9).Value = 32300
$ws.Cells.Item(24, 10).Value = 53800
$ws.Cells.Item(24, 11).Value = 84900

$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0

$ws.Cells.Item(26, 4).Value = -17700
$ws.Cells.Item(26, 5).Value = 31200
$ws.Cells.Item(26, 6).Value = 41400
$ws.Cells.Item(26, 7).Value = 26800
$ws.Cells.Item(26, 8).Value = 97700
$ws.Cells.Item(26, 9).Value = 79500
$ws.Cells.Item(26, 10).Value = 146500
$ws.Cells.Item(26, 11).Value = 230400

$ws.Cells.Item(27, 4).Value = -17100
$ws.Cells.Item(27, 5).Value = 30900
$ws.Cells.Item(27, 6).Value = 40800
$ws.Cells.Item(27, 7).Value = 26200
$ws.Cells.Item(27, 8).Value = 96800
$ws.Cells.Item(27, 9).Value = 79200
$ws.Cells.Item(27, 10).Value = 146200
$ws.Cells.Item(27, 11).Value = 230200

$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0

$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0

$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0

$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0

$ws.Cells.Item(32, 4).Value = 1700
$ws.Cells.Item(32, 5).Value = -800
$ws.Cells.Item(32, 6).Value = -1900
$ws.Cells.Item(32, 7).Value = 2200
$ws.Cells.Item(32, 8).Value = -14600
$ws.Cells.Item(32, 9).Value = -17100
$ws.Cells.Item(32, 10).Value = -24000
$ws.Cells.Item(32, 11).Value = -32400

$ws.Cells.Item(33, 4).Value = -17100
$ws.Cells.Item(33, 5).Value = 30900
$ws.Cells.Item(33, 6).Value = 40800
$ws.Cells.Item(33, 7).Value = 26200
$ws.Cells.Item(33, 8).Value = 96800
$ws.Cells.Item(33, 9).Value = 79200
$ws.Cells.Item(33, 10).Value = 146200
$ws.Cells.Item(33, 11).Value = 230200

$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0

$ws.Cells.Item(35, 4).Value = -17100
$ws.Cells.Item(35, 5).Value = 30900
$ws.Cells.Item(35, 6).Value = 40800
$ws.Cells.Item(35, 7).Value = 26200
$ws.Cells.Item(35, 8).Value = 96800
$ws.Cells.Item(35, 9).Value = 79200
$ws.Cells.Item(35, 10).Value = 146200
$ws.Cells.Item(35, 11).Value = 230200

$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43100
$ws.Cells.Item(38, 6).Value = 42735
$ws.Cells.Item(38, 7).Value = 42369
$ws.Cells.Item(38, 8).Value = 42004
$ws.Cells.Item(38, 9).Value = 41639
$ws.Cells.Item(38, 10).Value = 41274
$ws.Cells.Item(38, 11).Value = 40908

$ws.Cells.Item(41, 4).Value = 13800
$ws.Cells.Item(41, 5).Value = 16300
$ws.Cells.Item(41, 6).Value = 49100
$ws.Cells.Item(41, 7).Value = 15700
$ws.Cells.Item(41, 8).Value = 36700
$ws.Cells.Item(41, 9).Value = 16000
$ws.Cells.Item(41, 10).Value = 10500
$ws.Cells.Item(41, 11).Value = 3900

$ws.Cells.Item(42, 4).Value = 429300
$ws.Cells.Item(42, 5).Value = 229500
$ws.Cells.Item(42, 6).Value = 371400
$ws.Cells.Item(42, 7).Value = 271800
$ws.Cells.Item(42, 8).Value = 426500
$ws.Cells.Item(42, 9).Value = 417200
$ws.Cells.Item(42, 10).Value = 515700
$ws.Cells.Item(42, 11).Value = 327400

$ws.Cells.Item(43, 4).Value = 322800
$ws.Cells.Item(43, 5).Value = 356700
$ws.Cells.Item(43, 6).Value = 367500
$ws.Cells.Item(43, 7).Value = 362800
$ws.Cells.Item(43, 8).Value = 393200
$ws.Cells.Item(43, 9).Value = 480300
$ws.Cells.Item(43, 10).Value = 515400
$ws.Cells.Item(43, 11).Value = 521000

$ws.Cells.Item(44, 4).Value = 25900
$ws.Cells.Item(44, 5).Value = 19300
$ws.Cells.Item(44, 6).Value = 17700
$ws.Cells.Item(44, 7).Value = 23500
$ws.Cells.Item(44, 8).Value = 24600
$ws.Cells.Item(44, 9).Value = 23200
$ws.Cells.Item(44, 10).Value = 30900
$ws.Cells.Item(44, 11).Value = 47500

$ws.Cells.Item(45, 4).Value = 21200
$ws.Cells.Item(45, 5).Value = 22100
$ws.Cells.Item(45, 6).Value = 21800
$ws.Cells.Item(45, 7).Value = 11900
$ws.Cells.Item(45, 8).Value = 15700
$ws.Cells.Item(45, 9).Value = 14600
$ws.Cells.Item(45, 10).Value = 15200
$ws.Cells.Item(45, 11).Value = 31000

$ws.Cells.Item(46, 4).Value = 813100
$ws.Cells.Item(46, 5).Value = 644000
$ws.Cells.Item(46, 6).Value = 827400
$ws.Cells.Item(46, 7).Value = 685600
$ws.Cells.Item(46, 8).Value = 896700
$ws.Cells.Item(46, 9).Value = 951300
$ws.Cells.Item(46, 10).Value = 1087600
$ws.Cells.Item(46, 11).Value = 930700

$ws.Cells.Item(47, 4).Value = 235100
$ws.Cells.Item(47, 5).Value = 246900
$ws.Cells.Item(47, 6).Value = 219600
$ws.Cells.Item(47, 7).Value = 216600
$ws.Cells.Item(47, 8).Value = 227300
$ws.Cells.Item(47, 9).Value = 235600
$ws.Cells.Item(47, 10).Value = 336300
$ws.Cells.Item(47, 11).Value = 373500

$ws.Cells.Item(48, 4).Value = 455800
$ws.Cells.Item(48, 5).Value = 440900
$ws.Cells.Item(48, 6).Value = 457700
$ws.Cells.Item(48, 7).Value = 481400
$ws.Cells.Item(48, 8).Value = 506000
$ws.Cells.Item(48, 9).Value = 514600
$ws.Cells.Item(48, 10).Value = 573000
$ws.Cells.Item(48, 11).Value = 605600

$ws.Cells.Item(49, 4).Value = 358100
$ws.Cells.Item(49, 5).Value = 347600
$ws.Cells.Item(49, 6).Value = 333000
$ws.Cells.Item(49, 7).Value = 346000
$ws.Cells.Item(49, 8).Value = 362800
$ws.Cells.Item(49, 9).Value = 383200
$ws.Cells.Item(49, 10).Value = 416900
$ws.Cells.Item(49, 11).Value = 454700

$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0

$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0

$ws.Cells.Item(52, 4).Value = "NA"
$ws.Cells.Item(52, 5).Value = 0
$ws.Cells.Item(52, 6).Value = 300
$ws.Cells.Item(52, 7).Value = 2500
$ws.Cells.Item(52, 8).Value = 4700
$ws.Cells.Item(52, 9).Value = 6300
$ws.Cells.Item(52, 10).Value = 10500
$ws.Cells.Item(52, 11).Value = 25700

$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 0

$ws.Cells.Item(54, 4).Value = 1862000
$ws.Cells.Item(54, 5).Value = 1679400
$ws.Cells.Item(54, 6).Value = 1838000
$ws.Cells.Item(54, 7).Value = 1732100
$ws.Cells.Item(54, 8).Value = 1997500
$ws.Cells.Item(54, 9).Value = 2091000
$ws.Cells.Item(54, 10).Value = 2424300
$ws.Cells.Item(54, 11).Value = 2390200

$ws.Cells.Item(57, 4).Value = 79500
$ws.Cells.Item(57, 5).Value = 70600
$ws.Cells.Item(57, 6).Value = 56000
$ws.Cells.Item(57, 7).Value = 62600
$ws.Cells.Item(57, 8).Value = 121700
$ws.Cells.Item(57, 9).Value = 78900
$ws.Cells.Item(57, 10).Value = 112300
$ws.Cells.Item(57, 11).Value = 154700

$ws.Cells.Item(58, 4).Value = 171100
$ws.Cells.Item(58, 5).Value = 170500
$ws.Cells.Item(58, 6).Value = 238100
$ws.Cells.Item(58, 7).Value = 202500
$ws.Cells.Item(58, 8).Value = 301300
$ws.Cells.Item(58, 9).Value = 303500
$ws.Cells.Item(58, 10).Value = 311500
$ws.Cells.Item(58, 11).Value = 188300

$ws.Cells.Item(59, 4).Value = 212400
$ws.Cells.Item(59, 5).Value = 211900
$ws.Cells.Item(59, 6).Value = 237000
$ws.Cells.Item(59, 7).Value = 248000
$ws.Cells.Item(59, 8).Value = 242800
$ws.Cells.Item(59, 9).Value = 270400
$ws.Cells.Item(59, 10).Value = 323900
$ws.Cells.Item(59, 11).Value = 398000

$ws.Cells.Item(60, 4).Value = 463000
$ws.Cells.Item(60, 5).Value = 453000
$ws.Cells.Item(60, 6).Value = 531100
$ws.Cells.Item(60, 7).Value = 513200
$ws.Cells.Item(60, 8).Value = 665700
$ws.Cells.Item(60, 9).Value = 652800
$ws.Cells.Item(60, 10).Value = 747700
$ws.Cells.Item(60, 11).Value = 741100

$ws.Cells.Item(61, 4).Value = 895300
$ws.Cells.Item(61, 5).Value = 778600
$ws.Cells.Item(61, 6).Value = 790700
$ws.Cells.Item(61, 7).Value = 842600
$ws.Cells.Item(61, 8).Value = 978900
$ws.Cells.Item(61, 9).Value = 1196600
$ws.Cells.Item(61, 10).Value = 1483800
$ws.Cells.Item(61, 11).Value = 1528200

$ws.Cells.Item(62, 4).Value = 41100
$ws.Cells.Item(62, 5).Value = 50200
$ws.Cells.Item(62, 6).Value = 146500
$ws.Cells.Item(62, 7).Value = 49400
$ws.Cells.Item(62, 8).Value = 51600
$ws.Cells.Item(62, 9).Value = 45800
$ws.Cells.Item(62, 10).Value = 54900
$ws.Cells.Item(62, 11).Value = 68700

$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0

$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0

$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0

$ws.Cells.Item(66, 4).Value = 1399900
$ws.Cells.Item(66, 5).Value = 1282900
$ws.Cells.Item(66, 6).Value = 1473300
$ws.Cells.Item(66, 7).Value = 1409600
$ws.Cells.Item(66, 8).Value = 1700600
$ws.Cells.Item(66, 9).Value = 1896000
$ws.Cells.Item(66, 10).Value = 2286900
$ws.Cells.Item(66, 11).Value = 2339100

$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0

$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0

$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0

$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0

$ws.Cells.Item(72, 4).Value = 372200
$ws.Cells.Item(72, 5).Value = 396200
$ws.Cells.Item(72, 6).Value = 364700
$ws.Cells.Item(72, 7).Value = 322800
$ws.Cells.Item(72, 8).Value = 297400
$ws.Cells.Item(72, 9).Value = 198400
$ws.Cells.Item(72, 10).Value = 140400
$ws.Cells.Item(72, 11).Value = 50800

$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0

$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0

$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0

$ws.Cells.Item(76, 4).Value = 462100
$ws.Cells.Item(76, 5).Value = 396500
$ws.Cells.Item(76, 6).Value = 364700
$ws.Cells.Item(76, 7).Value = 322500
$ws.Cells.Item(76, 8).Value = 296900
$ws.Cells.Item(76, 9).Value = 195100
$ws.Cells.Item(76, 10).Value = 137400
$ws.Cells.Item(76, 11).Value = 51100

$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0

$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43100
$ws.Cells.Item(80, 6).Value = 42735
$ws.Cells.Item(80, 7).Value = 42369
$ws.Cells.Item(80, 8).Value = 42004
$ws.Cells.Item(80, 9).Value = 41639
$ws.Cells.Item(80, 10).Value = 41274
$ws.Cells.Item(80, 11).Value = 40908

$ws.Cells.Item(81, 4).Value = -17100
$ws.Cells.Item(81, 5).Value = 30900
$ws.Cells.Item(81, 6).Value = 40800
$ws.Cells.Item(81, 7).Value = 26200
$ws.Cells.Item(81, 8).Value = 96800
$ws.Cells.Item(81, 9).Value = 79200
$ws.Cells.Item(81, 10).Value = 146200
$ws.Cells.Item(81, 11).Value = 230200

$ws.Cells.Item(83, 4).Value = 161100
$ws.Cells.Item(83, 5).Value = 153100
$ws.Cells.Item(83, 6).Value = 147300
$ws.Cells.Item(83, 7).Value = 155100
$ws.Cells.Item(83, 8).Value = 168300
$ws.Cells.Item(83, 9).Value = 186500
$ws.Cells.Item(83, 10).Value = 211100
$ws.Cells.Item(83, 11).Value = 206100

$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0

$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0

$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0

$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0

$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0

$ws.Cells.Item(89, 4).Value = 212400
$ws.Cells.Item(89, 5).Value = 213500
$ws.Cells.Item(89, 6).Value = 215500
$ws.Cells.Item(89, 7).Value = 230700
$ws.Cells.Item(89, 8).Value = 429600
$ws.Cells.Item(89, 9).Value = 429300
$ws.Cells.Item(89, 10).Value = 452800
$ws.Cells.Item(89, 11).Value = 371000

$ws.Cells.Item(91, 4).Value = -98200
$ws.Cells.Item(91, 5).Value = -95500
$ws.Cells.Item(91, 6).Value = -81400
$ws.Cells.Item(91, 7).Value = -84100
$ws.Cells.Item(91, 8).Value = -79700
$ws.Cells.Item(91, 9).Value = -75900
$ws.Cells.Item(91, 10).Value = -126100
$ws.Cells.Item(91, 11).Value = -120700

$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0

$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0

$ws.Cells.Item(94, 4).Value = -174100
$ws.Cells.Item(94, 5).Value = -177700
$ws.Cells.Item(94, 6).Value = -100400
$ws.Cells.Item(94, 7).Value = -26500
$ws.Cells.Item(94, 8).Value = -96600
$ws.Cells.Item(94, 9).Value = -94900
$ws.Cells.Item(94, 10).Value = -195300
$ws.Cells.Item(94, 11).Value = -462600

$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = -300
$ws.Cells.Item(96, 6).Value = -300
$ws.Cells.Item(96, 7).Value = -300
$ws.Cells.Item(96, 8).Value = -1100
$ws.Cells.Item(96, 9).Value = -22300
$ws.Cells.Item(96, 10).Value = -107900
$ws.Cells.Item(96, 11).Value = -239700

$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0

$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0

$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0

$ws.Cells.Item(100, 4).Value = 148200
$ws.Cells.Item(100, 5).Value = -232600
$ws.Cells.Item(100, 6).Value = 17100
$ws.Cells.Item(100, 7).Value = -313400
$ws.Cells.Item(100, 8).Value = -305100
$ws.Cells.Item(100, 9).Value = -432900
$ws.Cells.Item(100, 10).Value = -121100
$ws.Cells.Item(100, 11).Value = 199700

$ws.Cells.Item(101, 4).Value = -300
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = -300
$ws.Cells.Item(101, 8).Value = "NA"
$ws.Cells.Item(101, 9).Value = "NA"
$ws.Cells.Item(101, 10).Value = "NA"
$ws.Cells.Item(101, 11).Value = "NA"

$ws.Cells.Item(102, 4).Value = 186200
$ws.Cells.Item(102, 5).Value = -196700
$ws.Cells.Item(102, 6).Value = 132200
$ws.Cells.Item(102, 7).Value = -109500
$ws.Cells.Item(102, 8).Value = 27900
$ws.Cells.Item(102, 9).Value = -98500
$ws.Cells.Item(102, 10).Value = 136300
$ws.Cells.Item(102, 11).Value = 108100

Write-Host "Edit complete"
